$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 156.85715
$ws.Range("I12").Value = 139.6
$ws.Range("K12").Value = 139.6
$ws.Range("M12").Value = 30.40000000000001
$ws.Range("H28").Value = 3573.4375
$ws.Range("J28").Value = 5074.6665
$ws.Range("L28").Value = 5074.6665
$ws.Range("N28").Value = -6044.6665
$ws.Range("H33").Value = 359.2857
$ws.Range("I33").Value = 239.33333
$ws.Range("K33").Value = 239.33333
$ws.Range("M33").Value = -10.33332999999999
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H64").Value = 3664.8
$ws.Range("I64").Value = 3872
$ws.Range("K64").Value = 3872
$ws.Range("M64").Value = -3624
$ws.Range("H67").Value = 3664.8
$ws.Range("I67").Value = 3872
$ws.Range("K67").Value = 3872
$ws.Range("M67").Value = -3014
$ws.Range("H112").Value = 2621.889
$ws.Range("J112").Value = 3171.2144
$ws.Range("L112").Value = 9513.643199999999
$ws.Range("N112").Value = -11729.6432
$ws.Range("H113").Value = 1249.5
$ws.Range("I113").Value = 1249.5
$ws.Range("K113").Value = 1249.5
$ws.Range("M113").Value = 2004.5
$ws.Range("H129").Value = 759.9524
$ws.Range("I129").Value = 403.7
$ws.Range("J129").Value = 871.28125
$ws.Range("K129").Value = 1211.1
$ws.Range("L129").Value = 2613.84375
$ws.Range("M129").Value = 3788.9
$ws.Range("N129").Value = -12613.84375
$ws.Range("H132").Value = 10425099
$ws.Range("I132").Value = 13896344
$ws.Range("J132").Value = 11364.75
$ws.Range("K132").Value = 41689032
$ws.Range("L132").Value = 34094.25
$ws.Range("M132").Value = -41686502
$ws.Range("N132").Value = -39154.25
$ws.Range("H135").Value = 31250814
$ws.Range("I135").Value = 613.86365
$ws.Range("K135").Value = 5524.77285
$ws.Range("M135").Value = -2989.77285

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 688.36584
$ws.Range("I2").Value = 460.42856
$ws.Range("J2").Value = 1179.3077
$ws.Range("K2").Value = 460.42856
$ws.Range("L2").Value = 1179.3077
$ws.Range("M2").Value = -347.42856
$ws.Range("N2").Value = -1405.3077
$ws.Range("H5").Value = 205.15384
$ws.Range("I5").Value = 151.85715
$ws.Range("J5").Value = 267.33334
$ws.Range("K5").Value = 151.85715
$ws.Range("L5").Value = 267.33334
$ws.Range("M5").Value = -39.85714999999999
$ws.Range("N5").Value = -491.33334
$ws.Range("H32").Value = 7923.8535
$ws.Range("I32").Value = 6648.597
$ws.Range("K32").Value = 6648.597
$ws.Range("M32").Value = -6361.597
$ws.Range("H45").Value = 1124
$ws.Range("I45").Value = 1063.7333
$ws.Range("J45").Value = 1350
$ws.Range("K45").Value = 1063.7333
$ws.Range("L45").Value = 1350
$ws.Range("M45").Value = -686.7333000000001
$ws.Range("N45").Value = -2104
$ws.Range("H102").Value = 8334508
$ws.Range("I102").Value = 11112124
$ws.Range("J102").Value = 1659.6
$ws.Range("K102").Value = 11112124
$ws.Range("L102").Value = 1659.6
$ws.Range("M102").Value = -11110502
$ws.Range("N102").Value = -4903.6
$ws.Range("H116").Value = 688.36584
$ws.Range("I116").Value = 460.42856
$ws.Range("J116").Value = 1179.3077
$ws.Range("K116").Value = 460.42856
$ws.Range("L116").Value = 1179.3077
$ws.Range("M116").Value = 1833.57144
$ws.Range("N116").Value = -5767.3077

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 688.36584
$ws.Range("I3").Value = 460.42856
$ws.Range("J3").Value = 1179.3077
$ws.Range("K3").Value = 460.42856
$ws.Range("L3").Value = 1179.3077
$ws.Range("M3").Value = -346.42856
$ws.Range("N3").Value = -1407.3077
$ws.Range("H4").Value = 205.15384
$ws.Range("I4").Value = 151.85715
$ws.Range("J4").Value = 267.33334
$ws.Range("K4").Value = 151.85715
$ws.Range("L4").Value = 267.33334
$ws.Range("M4").Value = -36.85714999999999
$ws.Range("N4").Value = -497.33334
$ws.Range("H22").Value = 466.66666
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227
$ws.Range("H105").Value = 200000900
$ws.Range("I105").Value = 333333660
$ws.Range("J105").Value = 1750
$ws.Range("K105").Value = 333333660
$ws.Range("L105").Value = 1750
$ws.Range("M105").Value = -333331913
$ws.Range("N105").Value = -5244

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 777.8570999999999
$ws.Range("I122").Value = 757.5
$ws.Range("K122").Value = 2272.5
$ws.Range("M122").Value = 177.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 104167500
$ws.Range("J129").Value = 27778886
$ws.Range("L129").Value = 83336658
$ws.Range("N129").Value = -83346658
$ws.Range("H131").Value = 18870610
$ws.Range("J131").Value = 3264.0952
$ws.Range("L131").Value = 9792.285600000001
$ws.Range("N131").Value = -19872.2856
$ws.Range("H136").Value = 2979.6667
$ws.Range("I136").Value = 2866.25
$ws.Range("J136").Value = 3109.2856
$ws.Range("K136").Value = 8598.75
$ws.Range("L136").Value = 9327.856800000001
$ws.Range("M136").Value = -3498.75
$ws.Range("N136").Value = -19527.8568

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5937.5
$ws.Range("I80").Value = 7866.6665
$ws.Range("J80").Value = 4780
$ws.Range("K80").Value = 7866.6665
$ws.Range("L80").Value = 4780
$ws.Range("M80").Value = -6868.6665
$ws.Range("N80").Value = -6776
$ws.Range("H83").Value = 5937.5
$ws.Range("I83").Value = 7866.6665
$ws.Range("J83").Value = 4780
$ws.Range("K83").Value = 39333.3325
$ws.Range("L83").Value = 23900
$ws.Range("M83").Value = -34341.3325
$ws.Range("N83").Value = -33884
$ws.Range("H102").Value = 1430.75
$ws.Range("I102").Value = 1392.0769
$ws.Range("K102").Value = 1392.0769
$ws.Range("M102").Value = 229.9231
$ws.Range("H113").Value = 1443.2916
$ws.Range("I113").Value = 1291.2142
$ws.Range("J113").Value = 1656.2
$ws.Range("K113").Value = 1291.2142
$ws.Range("L113").Value = 1656.2
$ws.Range("M113").Value = 878.7858000000001
$ws.Range("N113").Value = -5996.2
$ws.Range("H122").Value = 4810.722
$ws.Range("I122").Value = 5045.6924
$ws.Range("K122").Value = 15137.0772
$ws.Range("M122").Value = -12687.0772
$ws.Range("H126").Value = 1859.0435
$ws.Range("I126").Value = 1573.7333
$ws.Range("J126").Value = 2394
$ws.Range("K126").Value = 4721.199900000001
$ws.Range("L126").Value = 7182
$ws.Range("M126").Value = -2251.199900000001
$ws.Range("N126").Value = -12122
$ws.Range("H132").Value = 2689.6
$ws.Range("I132").Value = 2587.4546
$ws.Range("K132").Value = 7762.3638
$ws.Range("M132").Value = -5232.3638

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4521.615
$ws.Range("I46").Value = 745.25
$ws.Range("J46").Value = 6200
$ws.Range("K46").Value = 745.25
$ws.Range("L46").Value = 6200
$ws.Range("M46").Value = -557.25
$ws.Range("N46").Value = -6576
$ws.Range("H61").Value = 1298
$ws.Range("I61").Value = 1145.6666
$ws.Range("J61").Value = 1602.6666
$ws.Range("K61").Value = 1145.6666
$ws.Range("L61").Value = 1602.6666
$ws.Range("M61").Value = -943.6666
$ws.Range("N61").Value = -2006.6666
$ws.Range("H100").Value = 994.5263
$ws.Range("I100").Value = 873.06665
$ws.Range("K100").Value = 873.06665
$ws.Range("M100").Value = -332.06665
$ws.Range("H113").Value = 1298
$ws.Range("I113").Value = 1145.6666
$ws.Range("J113").Value = 1602.6666
$ws.Range("K113").Value = 1145.6666
$ws.Range("L113").Value = 1602.6666
$ws.Range("M113").Value = 1024.3334
$ws.Range("N113").Value = -5942.6666
$ws.Range("H122").Value = 31251788
$ws.Range("I122").Value = 62501576
$ws.Range("K122").Value = 187504728
$ws.Range("M122").Value = -187502278
$ws.Range("H136").Value = 2135.8667
$ws.Range("I136").Value = 1430.5
$ws.Range("K136").Value = 4291.5
$ws.Range("M136").Value = -1741.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H96").Value = 3598.75
$ws.Range("I96").Value = 3453.5454
$ws.Range("J96").Value = 3776.2222
$ws.Range("K96").Value = 3453.5454
$ws.Range("L96").Value = 3776.2222
$ws.Range("M96").Value = -2080.5454
$ws.Range("N96").Value = -6522.2222
$ws.Range("H107").Value = 442.45456
$ws.Range("I107").Value = 331.0909
$ws.Range("J107").Value = 553.8182
$ws.Range("K107").Value = 993.2727
$ws.Range("L107").Value = 1661.4546
$ws.Range("M107").Value = 926.7273
$ws.Range("N107").Value = -5501.4546
$ws.Range("H113").Value = 565.3
$ws.Range("I113").Value = 408.16666
$ws.Range("J113").Value = 801
$ws.Range("K113").Value = 1224.49998
$ws.Range("L113").Value = 2403
$ws.Range("M113").Value = 945.5000199999999
$ws.Range("N113").Value = -6743
$ws.Range("H122").Value = 10871256
$ws.Range("I122").Value = 11906494
$ws.Range("J122").Value = 1252
$ws.Range("K122").Value = 35719482
$ws.Range("L122").Value = 3756
$ws.Range("M122").Value = -35717032
$ws.Range("N122").Value = -8656
$ws.Range("H136").Value = 1442.7916
$ws.Range("I136").Value = 1226.0625
$ws.Range("J136").Value = 1876.25
$ws.Range("K136").Value = 3678.1875
$ws.Range("L136").Value = 5628.75
$ws.Range("M136").Value = -1128.1875
$ws.Range("N136").Value = -10728.75

